$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.037948750878991
$ws.Cells.Item(2, 4).Value2 = 1.038545000071734
$ws.Cells.Item(2, 5).Value2 = 1.04551580758189
$ws.Cells.Item(2, 6).Value2 = 1.053702120090383
$ws.Cells.Item(2, 9).Value2 = 1.033224874792498
$ws.Cells.Item(2, 10).Value2 = 1.043048846335555
$ws.Cells.Item(2, 11).Value2 = 1.041332728820576
$ws.Cells.Item(2, 12).Value2 = 1.048283829613246
$ws.Cells.Item(2, 13).Value2 = 1.056447359988555
$ws.Cells.Item(2, 14).Value2 = 1.018211673902787

$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.039128092039523
$ws.Cells.Item(3, 4).Value2 = 1.039577915622653
$ws.Cells.Item(3, 5).Value2 = 1.046601494685241
$ws.Cells.Item(3, 6).Value2 = 1.054960298436522
$ws.Cells.Item(3, 9).Value2 = 1.033391067949188
$ws.Cells.Item(3, 10).Value2 = 1.043871624532185
$ws.Cells.Item(3, 11).Value2 = 1.04217509102618
$ws.Cells.Item(3, 12).Value2 = 1.049180259058331
$ws.Cells.Item(3, 13).Value2 = 1.057517496391835
$ws.Cells.Item(3, 14).Value2 = 1.018490192386643

$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.039890906351399
$ws.Cells.Item(4, 4).Value2 = 1.04024627920079
$ws.Cells.Item(4, 5).Value2 = 1.047304125881551
$ws.Cells.Item(4, 6).Value2 = 1.055774835342065
$ws.Cells.Item(4, 9).Value2 = 1.033496934626027
$ws.Cells.Item(4, 10).Value2 = 1.044403203667407
$ws.Cells.Item(4, 11).Value2 = 1.042719552710937
$ws.Cells.Item(4, 12).Value2 = 1.049759835759399
$ws.Cells.Item(4, 13).Value2 = 1.05820979905978
$ws.Cells.Item(4, 14).Value2 = 1.018669994418009

$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.040211523712528
$ws.Cells.Item(5, 4).Value2 = 1.04052725957422
$ws.Cells.Item(5, 5).Value2 = 1.047599541462755
$ws.Cells.Item(5, 6).Value2 = 1.05611736649795
$ws.Cells.Item(5, 9).Value2 = 1.033541040870176
$ws.Cells.Item(5, 10).Value2 = 1.044626485602743
$ws.Cells.Item(5, 11).Value2 = 1.042948300787479
$ws.Cells.Item(5, 12).Value2 = 1.050003377404704
$ws.Cells.Item(5, 13).Value2 = 1.058500808998235
$ws.Cells.Item(5, 14).Value2 = 1.018745483284292

$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.040265352734026
$ws.Cells.Item(6, 4).Value2 = 1.040574437456651
$ws.Cells.Item(6, 5).Value2 = 1.047649144791216
$ws.Cells.Item(6, 6).Value2 = 1.056174884880314
$ws.Cells.Item(6, 9).Value2 = 1.033548423041081
$ws.Cells.Item(6, 10).Value2 = 1.044663964272774
$ws.Cells.Item(6, 11).Value2 = 1.042986700195867
$ws.Cells.Item(6, 12).Value2 = 1.050044262552348
$ws.Cells.Item(6, 13).Value2 = 1.058549668880203
$ws.Cells.Item(6, 14).Value2 = 1.018758152343519

$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.039895190728403
$ws.Cells.Item(7, 4).Value2 = 1.040250033672062
$ws.Cells.Item(7, 5).Value2 = 1.047308073122227
$ws.Cells.Item(7, 6).Value2 = 1.055779411866137
$ws.Cells.Item(7, 9).Value2 = 1.033497525547931
$ws.Cells.Item(7, 10).Value2 = 1.044406187932603
$ws.Cells.Item(7, 11).Value2 = 1.042722609817328
$ws.Cells.Item(7, 12).Value2 = 1.049763090415373
$ws.Cells.Item(7, 13).Value2 = 1.058213687681425
$ws.Cells.Item(7, 14).Value2 = 1.018671003496289

$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.038347377490157
$ws.Cells.Item(8, 4).Value2 = 1.038894079565803
$ws.Cells.Item(8, 5).Value2 = 1.045882696532272
$ws.Cells.Item(8, 6).Value2 = 1.054127242989604
$ws.Cells.Item(8, 9).Value2 = 1.03328138664405
$ws.Cells.Item(8, 10).Value2 = 1.043327076971219
$ws.Cells.Item(8, 11).Value2 = 1.041617534372154
$ws.Cells.Item(8, 12).Value2 = 1.048586880616547
$ws.Cells.Item(8, 13).Value2 = 1.056809048740204
$ws.Cells.Item(8, 14).Value2 = 1.018305887194619

$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.035617574848984
$ws.Cells.Item(9, 4).Value2 = 1.036504665983075
$ws.Cells.Item(9, 5).Value2 = 1.043371870536438
$ws.Cells.Item(9, 6).Value2 = 1.051218999959192
$ws.Cells.Item(9, 9).Value2 = 1.03288771951435
$ws.Cells.Item(9, 10).Value2 = 1.041419280179384
$ws.Cells.Item(9, 11).Value2 = 1.039665598353855
$ws.Cells.Item(9, 12).Value2 = 1.046510586377913
$ws.Cells.Item(9, 13).Value2 = 1.054332717773719
$ws.Cells.Item(9, 14).Value2 = 1.017659295767725

$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.03379601575523
$ws.Cells.Item(10, 4).Value2 = 1.034911638659489
$ws.Cells.Item(10, 5).Value2 = 1.041698509644435
$ws.Cells.Item(10, 6).Value2 = 1.049282160393913
$ws.Cells.Item(10, 9).Value2 = 1.032616663060303
$ws.Cells.Item(10, 10).Value2 = 1.040143150421053
$ws.Cells.Item(10, 11).Value2 = 1.038361130332835
$ws.Cells.Item(10, 12).Value2 = 1.045123870908656
$ws.Cells.Item(10, 13).Value2 = 1.052680965061529
$ws.Cells.Item(10, 14).Value2 = 1.017226064345152

$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.033006835317778
$ws.Cells.Item(11, 4).Value2 = 1.034221807653594
$ws.Cells.Item(11, 5).Value2 = 1.040974035389579
$ws.Cells.Item(11, 6).Value2 = 1.048443938792777
$ws.Cells.Item(11, 9).Value2 = 1.032497248744063
$ws.Cells.Item(11, 10).Value2 = 1.039589547968832
$ws.Cells.Item(11, 11).Value2 = 1.037795515805884
$ws.Cells.Item(11, 12).Value2 = 1.044522796786003
$ws.Cells.Item(11, 13).Value2 = 1.051965515437222
$ws.Cells.Item(11, 14).Value2 = 1.017037952470461

$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.032713631270596
$ws.Cells.Item(12, 4).Value2 = 1.033965566568014
$ws.Cells.Item(12, 5).Value2 = 1.0407049473676
$ws.Cells.Item(12, 6).Value2 = 1.048132650615182
$ws.Cells.Item(12, 9).Value2 = 1.032452585521193
$ws.Cells.Item(12, 10).Value2 = 1.039383759466775
$ws.Cells.Item(12, 11).Value2 = 1.037585304393071
$ws.Cells.Item(12, 12).Value2 = 1.044299437112749
$ws.Cells.Item(12, 13).Value2 = 1.051699729483266
$ws.Cells.Item(12, 14).Value2 = 1.016968000942216

$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.032776527614489
$ws.Cells.Item(13, 4).Value2 = 1.034020531492545
$ws.Cells.Item(13, 5).Value2 = 1.040762667051309
$ws.Cells.Item(13, 6).Value2 = 1.048199420090576
$ws.Cells.Item(13, 9).Value2 = 1.032462179848805
$ws.Cells.Item(13, 10).Value2 = 1.039427908872928
$ws.Cells.Item(13, 11).Value2 = 1.037630400766937
$ws.Cells.Item(13, 12).Value2 = 1.044347352793979
$ws.Cells.Item(13, 13).Value2 = 1.051756743132067
$ws.Cells.Item(13, 14).Value2 = 1.016983009335303

$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.032982600365254
$ws.Cells.Item(14, 4).Value2 = 1.034200626837909
$ws.Cells.Item(14, 5).Value2 = 1.040951792212566
$ws.Cells.Item(14, 6).Value2 = 1.048418206306943
$ws.Cells.Item(14, 9).Value2 = 1.032493563140623
$ws.Cells.Item(14, 10).Value2 = 1.039572540614253
$ws.Cells.Item(14, 11).Value2 = 1.037778142058053
$ws.Cells.Item(14, 12).Value2 = 1.044504335730219
$ws.Cells.Item(14, 13).Value2 = 1.051943546233127
$ws.Cells.Item(14, 14).Value2 = 1.017032171856481

$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.033109559605441
$ws.Cells.Item(15, 4).Value2 = 1.034311588543296
$ws.Cells.Item(15, 5).Value2 = 1.041068320277209
$ws.Cells.Item(15, 6).Value2 = 1.048553016226147
$ws.Cells.Item(15, 9).Value2 = 1.032512858679808
$ws.Cells.Item(15, 10).Value2 = 1.039661632309613
$ws.Cells.Item(15, 11).Value2 = 1.037869154807228
$ws.Cells.Item(15, 12).Value2 = 1.04460104560718
$ws.Cells.Item(15, 13).Value2 = 1.052058636956849
$ws.Cells.Item(15, 14).Value2 = 1.017062452110063

$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.033848381745069
$ws.Cells.Item(16, 4).Value2 = 1.034957419535533
$ws.Cells.Item(16, 5).Value2 = 1.041746592649029
$ws.Cells.Item(16, 6).Value2 = 1.049337799569281
$ws.Cells.Item(16, 9).Value2 = 1.032624545092477
$ws.Cells.Item(16, 10).Value2 = 1.040179869394343
$ws.Cells.Item(16, 11).Value2 = 1.038398651961064
$ws.Cells.Item(16, 12).Value2 = 1.045163749108196
$ws.Cells.Item(16, 13).Value2 = 1.052728442120176
$ws.Cells.Item(16, 14).Value2 = 1.017238537732822

$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.034311707793131
$ws.Cells.Item(17, 4).Value2 = 1.03536252080642
$ws.Cells.Item(17, 5).Value2 = 1.042172081352698
$ws.Cells.Item(17, 6).Value2 = 1.049830190637231
$ws.Cells.Item(17, 9).Value2 = 1.032694055376301
$ws.Cells.Item(17, 10).Value2 = 1.040504669275512
$ws.Cells.Item(17, 11).Value2 = 1.038730584257634
$ws.Cells.Item(17, 12).Value2 = 1.045516552026092
$ws.Cells.Item(17, 13).Value2 = 1.05314853056242
$ws.Cells.Item(17, 14).Value2 = 1.017348852208055

$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.034581916377659
$ws.Cells.Item(18, 4).Value2 = 1.035598805778142
$ws.Cells.Item(18, 5).Value2 = 1.04242027165491
$ws.Cells.Item(18, 6).Value2 = 1.050117436859445
$ws.Cells.Item(18, 9).Value2 = 1.032734402240626
$ws.Cells.Item(18, 10).Value2 = 1.04069402026685
$ws.Cells.Item(18, 11).Value2 = 1.038924120411396
$ws.Cells.Item(18, 12).Value2 = 1.045722276537175
$ws.Cells.Item(18, 13).Value2 = 1.053393538979492
$ws.Cells.Item(18, 14).Value2 = 1.017413146642235

$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.034674043521554
$ws.Cells.Item(19, 4).Value2 = 1.035679372247004
$ws.Cells.Item(19, 5).Value2 = 1.042504899860693
$ws.Cells.Item(19, 6).Value2 = 1.05021538772084
$ws.Cells.Item(19, 9).Value2 = 1.032748126016986
$ws.Cells.Item(19, 10).Value2 = 1.040758567256495
$ws.Cells.Item(19, 11).Value2 = 1.038990098662594
$ws.Cells.Item(19, 12).Value2 = 1.045792413230823
$ws.Cells.Item(19, 13).Value2 = 1.053477076823039
$ws.Cells.Item(19, 14).Value2 = 1.017435060888456

$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.034262001619195
$ws.Cells.Item(20, 4).Value2 = 1.035319057670786
$ws.Cells.Item(20, 5).Value2 = 1.042126429424244
$ws.Cells.Item(20, 6).Value2 = 1.049777357318773
$ws.Cells.Item(20, 9).Value2 = 1.032686617990783
$ws.Cells.Item(20, 10).Value2 = 1.040469831600937
$ws.Cells.Item(20, 11).Value2 = 1.038694978761746
$ws.Cells.Item(20, 12).Value2 = 1.045478705772299
$ws.Cells.Item(20, 13).Value2 = 1.053103461342887
$ws.Cells.Item(20, 14).Value2 = 1.017337021696849

$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.032921918951669
$ws.Cells.Item(21, 4).Value2 = 1.034147593453907
$ws.Cells.Item(21, 5).Value2 = 1.040896099191603
$ws.Cells.Item(21, 6).Value2 = 1.048353777449766
$ws.Cells.Item(21, 9).Value2 = 1.032484330031376
$ws.Cells.Item(21, 10).Value2 = 1.039529954484447
$ws.Cells.Item(21, 11).Value2 = 1.037734639168788
$ws.Cells.Item(21, 12).Value2 = 1.044458110772665
$ws.Cells.Item(21, 13).Value2 = 1.051888538385311
$ws.Cells.Item(21, 14).Value2 = 1.017017696886102

$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.032078965356777
$ws.Cells.Item(22, 4).Value2 = 1.033411005536035
$ws.Cells.Item(22, 5).Value2 = 1.040122621911951
$ws.Cells.Item(22, 6).Value2 = 1.047459089383741
$ws.Cells.Item(22, 9).Value2 = 1.032355364783802
$ws.Cells.Item(22, 10).Value2 = 1.038938113848072
$ws.Cells.Item(22, 11).Value2 = 1.037130158104257
$ws.Cells.Item(22, 12).Value2 = 1.04381587781588
$ws.Cells.Item(22, 13).Value2 = 1.051124458403782
$ws.Cells.Item(22, 14).Value2 = 1.016816470938383

$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.03252586874883
$ws.Cells.Item(23, 4).Value2 = 1.033801489091396
$ws.Cells.Item(23, 5).Value2 = 1.040532649708941
$ws.Cells.Item(23, 6).Value2 = 1.047933345714698
$ws.Cells.Item(23, 9).Value2 = 1.032423900369504
$ws.Cells.Item(23, 10).Value2 = 1.039251945717911
$ws.Cells.Item(23, 11).Value2 = 1.037450669543289
$ws.Cells.Item(23, 12).Value2 = 1.044156389679589
$ws.Cells.Item(23, 13).Value2 = 1.051529532066057
$ws.Cells.Item(23, 14).Value2 = 1.016923187735433

$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.03428446182839
$ws.Cells.Item(24, 4).Value2 = 1.035338696802629
$ws.Cells.Item(24, 5).Value2 = 1.04214705753505
$ws.Cells.Item(24, 6).Value2 = 1.049801230288915
$ws.Cells.Item(24, 9).Value2 = 1.032689979234867
$ws.Cells.Item(24, 10).Value2 = 1.040485573553362
$ws.Cells.Item(24, 11).Value2 = 1.038711067582505
$ws.Cells.Item(24, 12).Value2 = 1.04549580704913
$ws.Cells.Item(24, 13).Value2 = 1.053123826251473
$ws.Cells.Item(24, 14).Value2 = 1.01734236755057

$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.036323584520675
$ws.Cells.Item(25, 4).Value2 = 1.037122397815831
$ws.Cells.Item(25, 5).Value2 = 1.04402088250408
$ws.Cells.Item(25, 6).Value2 = 1.051970493417908
$ws.Cells.Item(25, 9).Value2 = 1.032991009483172
$ws.Cells.Item(25, 10).Value2 = 1.041913239271943
$ws.Cells.Item(25, 11).Value2 = 1.040170776970207
$ws.Cells.Item(25, 12).Value2 = 1.047047798057645
$ws.Cells.Item(25, 13).Value2 = 1.054973055807876
$ws.Cells.Item(25, 14).Value2 = 1.01782683679626
